# XTXStarterKit-dev/scores.xlsx — "made submissions, awaiting blend10 submission"
#
# Rows 6-8 (trial 4,5,6 = blend 4,5 using lasso / lasso_drop / default) are
# marked as abandoned (strike-through, filled E:G, red "Error: Timed out"
# note in H) and re-run as trials 7-10 in rows 9-12 with the submission
# status tracked in columns H/K. Selection moves to E10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 6-8: mark the original blend-4,5 attempts as struck-through /
# timed out. B/C/D already hold values - we only change their
# formatting - then add new formatted-but-empty E:G cells and an H
# cell carrying the "Error: Timed out" note.
# ---------------------------------------------------------------------
foreach ($r in 6, 7, 8) {
    $ws.Range("B$r").Font.Strikethrough = $true
    $ws.Range("C$r").Font.Strikethrough = $true
    $ws.Range("D$r").Font.Strikethrough = $true
    $ws.Range("E$r").Font.Strikethrough = $true
    $ws.Range("F$r").Font.Strikethrough = $true
    $ws.Range("G$r").Font.Strikethrough = $true
}

# New red "Error: Timed out" note in column H for rows 6-8 (first write
# establishes the shared string + the new non-italic red font).
$ws.Range("H6").Value = "Error: Timed out"
$ws.Range("H6").Font.Italic = $false
$ws.Range("H6").Font.Color = 255
$ws.Range("H7").Value = "Error: Timed out"
$ws.Range("H7").Font.Italic = $false
$ws.Range("H7").Font.Color = 255
$ws.Range("H8").Value = "Error: Timed out"
$ws.Range("H8").Font.Italic = $false
$ws.Range("H8").Font.Color = 255

# ---------------------------------------------------------------------
# Rows 9-12: the re-submitted trials (7-10). Values/labels are written
# in an order chosen to reproduce the canonical shared-string table
# order (Error: Timed out / Runs..42s,500s / Runs..42s / Submitted /
# blend..lasso / blend..lasso_drop / Awaiting submission); formatting
# is applied afterwards.
# ---------------------------------------------------------------------
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "lasso"
$ws.Range("C9").Value = 0.050599999999999999
$ws.Range("D9").Value = -0.36285563865669301
$ws.Range("H9").Value = "Runs 100 in 42s, 500 in 500s"

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "lasso_drop"
$ws.Range("C10").Value = 0.058599999999999999
$ws.Range("D10").Value = -2.06886488756439
$ws.Range("H10").Value = "Runs 100 in 42s"

$ws.Range("K9").Value = "Submitted"

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "blend 4,5 with sigmoid_drop using lasso"
$ws.Range("D12").Value = 0.046677494892257802

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "blend 4,5 with sigmoid_drop using lasso_drop"
$ws.Range("D11").Value = -0.64641597628391601
$ws.Range("H11").Value = "Runs 100 in 42s"
$ws.Range("K11").Value = "Submitted"

$ws.Range("K10").Value = "Submitted"
$ws.Range("K12").Value = "Awaiting submission"

$ws.Range("A13").Value = 11

# Formatting for rows 9-12: B/C share one "plain" style, D carries the
# 0.0000 number format, C11/C12 keep the existing yellow highlight fill.
foreach ($r in 9, 10, 11, 12) {
    $ws.Range("B$r").Font.ThemeColor = 1
    $ws.Range("D$r").Font.ThemeColor = 1
    $ws.Range("D$r").NumberFormat = "0.0000"
}
$ws.Range("C9").Font.ThemeColor = 1
$ws.Range("C10").Font.ThemeColor = 1
$ws.Range("C11").Interior.Color = 65535
$ws.Range("C12").Interior.Color = 65535

# ---------------------------------------------------------------------
# Selection moves to E10
# ---------------------------------------------------------------------
$ws.Range("E10").Select()

Write-Host "edit complete"
